# Generate Report for Handback
#
# The handback files for both locales have now been received and are in
# sync with en-US, so the generated report is refreshed:
#  - Status moves from "Ready for handoff" to
#    "Handed back: in sync with en-US" (Overview sheet's per-language
#    status columns, and each language sheet's own Status column).
#  - "Latest Handback DateTime" is refreshed to the time the (now
#    up-to-date) handback file was generated, for each language.
#  - The previous "version not latest" Error Detail is cleared now that
#    the handback is in sync with the source.
#  - Columns that show the (now longer/shorter) Status and Error Detail
#    text are resized to fit the new content.

$wb = $excel.ActiveWorkbook

$statusNew = "Handed back: in sync with en-US"

# --- Overview sheet: per-language status cells ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusNew
$overview.Range("F2").Value = $statusNew
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusNew
$zhcn.Range("K2").Value = "2016-08-18 02:43:17"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.1
$zhcn.Columns.Item(16).ColumnWidth = 12.8

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusNew
$dede.Range("K2").Value = "2016-08-18 02:43:24"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.1
$dede.Columns.Item(16).ColumnWidth = 12.8
